$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Propagate the formatting used in row 2 (header-like style on column A,
# date style on column H) down to the two new rows, so the new cells
# reuse the existing style indices instead of creating new ones.
$ws.Range("A2").Copy()
$ws.Range("A3:A4").PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Range("H3:H4").PasteSpecial(-4122)

# --- Row 3 ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 28
$ws.Range("C3").Value = 13
$ws.Range("D3").Value = "s,vnds"
$ws.Range("E3").Value = "sdfnvkj"
$ws.Range("F3").Value = "djf nk"
$ws.Range("G3").Value = "Published"
$ws.Range("H3").Value = 44927

# --- Row 4 ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 29
$ws.Range("C4").Value = 17
$ws.Range("D4").Value = "New Patent"

# E4 ("5156") must stay a text value, not be auto-converted to a number.
# Build it as text in a scratch cell (apostrophe forces text), then copy
# only the value over to E4, leaving no stray formatting behind.
$ws.Range("Z1").Value = "'5156"
$ws.Range("Z1").Copy()
$ws.Range("E4").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("F4").Value = "askjcndskv"
$ws.Range("G4").Value = "Published"
$ws.Range("H4").Value = 44927
